$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range('D2').Value = '67.150.95'
$ws.Range('E2').Value = '  -1.97%  '
$ws.Range('D3').Value = '3.585.05'
$ws.Range('E3').Value = '  -3.24%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '575.71'
$ws.Range('E5').Value = '  -6.73%  '
$ws.Range('D6').Value = '191.81'
$ws.Range('E6').Value = '  -2.44%  '
$ws.Range('D7').Value = '3.580.30'
$ws.Range('E7').Value = '  -3.26%  '
$ws.Range('E8').Value = '  -2.68%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = '0.679'
$ws.Range('E10').Value = '  -6.35%  '
$ws.Range('E11').Value = '  -5.76%  '
$ws.Range('D12').Value = '56.06'
$ws.Range('E12').Value = '  -7.34%  '
$ws.Range('E13').Value = '  -6.23%  '
$ws.Range('D14').Value = '9.88'
$ws.Range('E14').Value = '  -5.31%  '
$ws.Range('D15').Value = '4.160.64'
$ws.Range('E15').Value = '  -3.28%  '
$ws.Range('D16').Value = '3.586.87'
$ws.Range('E16').Value = '  -3.20%  '
$ws.Range('D17').Value = '0.125'
$ws.Range('E17').Value = '  -1.43%  '
$ws.Range('D18').Value = '18.39'
$ws.Range('E18').Value = '  -5.25%  '
$ws.Range('D19').Value = '67.135.52'
$ws.Range('E19').Value = '  -1.82%  '
$ws.Range('D20').Value = '12.18'
$ws.Range('E20').Value = '  -5.32%  '
$ws.Range('D21').Value = '1.06'
$ws.Range('E21').Value = '  -7.10%  '
$ws.Range('D22').Value = '396.67'
$ws.Range('E22').Value = '  -2.93%  '
$ws.Range('D23').Value = '4.20'
$ws.Range('E23').Value = '  -9.90%  '
$ws.Range('D24').Value = '85.92'
$ws.Range('E24').Value = '  -4.45%  '
$ws.Range('D25').Value = '11.27'
$ws.Range('E25').Value = '  -1.75%  '
$ws.Range('D26').Value = '2.94'
$ws.Range('E26').Value = '  -4.31%  '
$ws.Range('D27').Value = '12.49'
$ws.Range('E27').Value = '  -4.35%  '
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('D29').Value = '3.61'
$ws.Range('E29').Value = '  -4.51%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = '7.75'
$ws.Range('E30').Value = '  +1.80%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '8.97'
$ws.Range('E31').Value = '  -6.97%  '
$ws.Range('D32').Value = '31.18'
$ws.Range('E32').Value = '  -4.63%  '
$ws.Range('D33').Value = '633.09'
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('D34').Value = '12.20'
$ws.Range('E35').Value = '  -6.37%  '
$ws.Range('D36').Value = '63.63'
$ws.Range('E36').Value = '  -5.93%  '
$ws.Range('D37').Value = '42.24'
$ws.Range('E37').Value = '  -9.22%  '
$ws.Range('D38').Value = '0.402'
$ws.Range('E38').Value = '  -2.48%  '
$ws.Range('E39').Value = '  +0.46%  '
$ws.Range('D40').Value = '0.0₃0766'
$ws.Range('E40').Value = '  -6.65%  '
$ws.Range('D41').Value = '0.134'
$ws.Range('E41').Value = '  -3.88%  '
$ws.Range('D42').Value = '3.137.96'
$ws.Range('E42').Value = '  +7.63%  '
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('D44').Value = '2.71'
$ws.Range('E44').Value = '  +3.46%  '
$ws.Range('D45').Value = '2.96'
$ws.Range('E45').Value = '  -2.75%  '
$ws.Range('E46').Value = '  -6.30%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '3.11'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '0.131'
$ws.Range('E48').Value = '  -6.67%  '
$ws.Range('D49').Value = '140.16'
$ws.Range('E49').Value = '  -4.54%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value = '2.58'
$ws.Range('E50').Value = '  +0.93%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '8.55'
$ws.Range('E51').Value = '  -9.19%  '
